# Log file updated, with links of Post45
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 55

# Shared-string table order mirrors the order these are first written:
# dev.to link (F), then title (C), then hashnode link (E).
$ws.Range("B$row").Value = 45
$ws.Range("F$row").Value = "https://dev.to/rahulmishra05/question-on-binary-semaphore-operating-system-m03-p07-3alp"
$ws.Range("C$row").Value = "Question on Binary Semaphore | Operating System - M03 P07"
$ws.Range("D$row").Value = 44170
$ws.Range("E$row").Value = "https://programmingport.hashnode.dev/question-on-binary-semaphore-or-operating-system-m03-p07"

# Match formatting used by the rest of the table (date format, hyperlink style)
$ws.Range("D$row").NumberFormat = "m/d/yy"
$ws.Range("E$row").Style = "Hyperlink"
$ws.Range("F$row").Style = "Hyperlink"

# Resize the table to include the newly added row
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("B10:F$row"))

$ws.Range("E$row").Select()
